# Apply odds updates to the FlashScore weekly games workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("G2").Value = 2.05
$ws.Range("I2").Value = 3.9
$ws.Range("X2").Value = 8.5
$ws.Range("AH2").Value = 19
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 51
$ws.Range("AR2").Value = 67

# Row 4 changes
$ws.Range("G4").Value = 3.1
$ws.Range("I4").Value = 2.35
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 3.2
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 15
$ws.Range("AJ4").Value = 21
$ws.Range("AT4").Value = 2.38

# Row 5 changes
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62

# Row 7 changes
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
